$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 250005500
$ws.Range("J94").Value = 333340000
$ws.Range("L94").Value = 333340000
$ws.Range("N94").Value = -333340902

$ws.Range("H132").Value = 9264922
$ws.Range("I132").Value = 9264922
$ws.Range("K132").Value = 27794766
$ws.Range("M132").Value = -27792236

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 2460
$ws.Range("I14").Value = 500
$ws.Range("J14").Value = 2950
$ws.Range("K14").Value = 500
$ws.Range("L14").Value = 2950
$ws.Range("M14").Value = -325
$ws.Range("N14").Value = -3300

$ws.Range("H132").Value = 3465.6428
$ws.Range("I132").Value = 3614.2856
$ws.Range("J132").Value = 3019.7144
$ws.Range("K132").Value = 10842.8568
$ws.Range("L132").Value = 9059.143199999999
$ws.Range("M132").Value = -8312.856800000001
$ws.Range("N132").Value = -14119.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 220960
$ws.Range("I86").Value = 275700
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 275700
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -274577
$ws.Range("N86").Value = -4246

$ws.Range("H89").Value = 220960
$ws.Range("I89").Value = 275700
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 1378500
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -1372884
$ws.Range("N89").Value = -21232

$ws.Range("H134").Value = 2873.6875
$ws.Range("I134").Value = 3020.3845
$ws.Range("J134").Value = 2238
$ws.Range("K134").Value = 9061.1535
$ws.Range("L134").Value = 6714
$ws.Range("M134").Value = -6526.1535
$ws.Range("N134").Value = -11784

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 11113231
$ws.Range("I62").Value = 27779028
$ws.Range("K62").Value = 27779028
$ws.Range("M62").Value = -27778404

$ws.Range("H65").Value = 11113231
$ws.Range("I65").Value = 27779028
$ws.Range("K65").Value = 138895140
$ws.Range("M65").Value = -138892020

$ws.Range("H134").Value = 955
$ws.Range("I134").Value = 839.3333
$ws.Range("K134").Value = 2517.9999
$ws.Range("M134").Value = 17.0001000000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3296.5
$ws.Range("I80").Value = 4980
$ws.Range("J80").Value = 2735.3333
$ws.Range("K80").Value = 4980
$ws.Range("L80").Value = 2735.3333
$ws.Range("M80").Value = -3982
$ws.Range("N80").Value = -4731.3333

$ws.Range("H83").Value = 3296.5
$ws.Range("I83").Value = 4980
$ws.Range("J83").Value = 2735.3333
$ws.Range("K83").Value = 24900
$ws.Range("L83").Value = 13676.6665
$ws.Range("M83").Value = -19908
$ws.Range("N83").Value = -23660.6665

$ws.Range("H132").Value = 2180.88
$ws.Range("I132").Value = 2118.5
$ws.Range("K132").Value = 6355.5
$ws.Range("M132").Value = -3825.5

$ws.Range("H134").Value = 26854.166
$ws.Range("J134").Value = 26854.166
$ws.Range("L134").Value = 80562.49800000001
$ws.Range("N134").Value = -85632.49800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1958.4667
$ws.Range("I7").Value = 1457.3158
$ws.Range("J7").Value = 2824.0908
$ws.Range("K7").Value = 1457.3158
$ws.Range("L7").Value = 2824.0908
$ws.Range("M7").Value = -1345.3158
$ws.Range("N7").Value = -3048.0908

$ws.Range("H22").Value = 592.2143
$ws.Range("I22").Value = 395.5
$ws.Range("K22").Value = 395.5
$ws.Range("M22").Value = -100.5

$ws.Range("H27").Value = 592.2143
$ws.Range("I27").Value = 395.5
$ws.Range("K27").Value = 395.5
$ws.Range("M27").Value = -288.5

$ws.Range("H46").Value = 389891.8
$ws.Range("I46").Value = 628.5714
$ws.Range("J46").Value = 533304.5600000001
$ws.Range("K46").Value = 628.5714
$ws.Range("L46").Value = 533304.5600000001
$ws.Range("M46").Value = -440.5714
$ws.Range("N46").Value = -533680.5600000001

$ws.Range("H68").Value = 4572.3335
$ws.Range("I68").Value = 2000.5
$ws.Range("K68").Value = 2000.5
$ws.Range("M68").Value = -1251.5

$ws.Range("H71").Value = 4572.3335
$ws.Range("I71").Value = 2000.5
$ws.Range("K71").Value = 10002.5
$ws.Range("M71").Value = -6258.5

$ws.Range("H82").Value = 1217.5
$ws.Range("I82").Value = 1267.75
$ws.Range("J82").Value = 1066.75
$ws.Range("K82").Value = 1267.75
$ws.Range("L82").Value = 1066.75
$ws.Range("M82").Value = -906.75
$ws.Range("N82").Value = -1788.75

$ws.Range("H85").Value = 1217.5
$ws.Range("I85").Value = 1267.75
$ws.Range("J85").Value = 1066.75
$ws.Range("K85").Value = 1267.75
$ws.Range("L85").Value = 1066.75
$ws.Range("M85").Value = -19.75
$ws.Range("N85").Value = -3562.75

$ws.Range("H122").Value = 4753.3335
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 4753.3335
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 14260.0005
$ws.Range("M122").Value = $null
$ws.Range("N122").Value = -19160.0005

$ws.Range("H126").Value = 1958.4667
$ws.Range("I126").Value = 1457.3158
$ws.Range("J126").Value = 2824.0908
$ws.Range("K126").Value = 4371.9474
$ws.Range("L126").Value = 8472.2724
$ws.Range("M126").Value = -1901.9474
$ws.Range("N126").Value = -13412.2724

$ws.Range("H132").Value = 3580.4666
$ws.Range("I132").Value = 4026.8948
$ws.Range("J132").Value = 2809.3635
$ws.Range("K132").Value = 12080.6844
$ws.Range("L132").Value = 8428.0905
$ws.Range("M132").Value = -9550.6844
$ws.Range("N132").Value = -13488.0905

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 24733
$ws.Range("I42").Value = 15000
$ws.Range("J42").Value = 26679.6
$ws.Range("K42").Value = 15000
$ws.Range("L42").Value = 26679.6
$ws.Range("M42").Value = -14622
$ws.Range("N42").Value = -27435.6

$ws.Range("H62").Value = 5002460
$ws.Range("I62").Value = 25001250
$ws.Range("J62").Value = 2762.5
$ws.Range("K62").Value = 25001250
$ws.Range("L62").Value = 2762.5
$ws.Range("M62").Value = -25000626
$ws.Range("N62").Value = -4010.5

$ws.Range("H65").Value = 5002460
$ws.Range("I65").Value = 25001250
$ws.Range("J65").Value = 2762.5
$ws.Range("K65").Value = 125006250
$ws.Range("L65").Value = 13812.5
$ws.Range("M65").Value = -125003130
$ws.Range("N65").Value = -20052.5

$ws.Range("H126").Value = 1033.9259
$ws.Range("I126").Value = 996.7826
$ws.Range("K126").Value = 2990.3478
$ws.Range("M126").Value = -520.3478
